$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 89) to the "jogos" table, mirroring the
# existing rows: Nome do Jogo | Status | Plataforma | Objetivo
$row = 89

# Column A holds a numeric-looking id ("888"). Several existing rows use
# plain numeric-looking ids stored as text (e.g. "1", "2", "666", "777"),
# so force a text number format before assigning the value to keep it text
# instead of Excel auto-coercing it to a number.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "888"

$ws.Cells.Item($row, 2).Value = "Incompleto"
$ws.Cells.Item($row, 3).Value = "PS3"
$ws.Cells.Item($row, 4).Value = "Zerar"
